$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append a new appointment record in row 3 ---------------------------
# Plain string literals like "False" / "None" would otherwise be
# auto-typed by the smart cell-input parser (e.g. "False" -> boolean
# FALSE). Route every value through a formula that evaluates to the
# literal text, then "Paste Special -> Values" over itself so the stored
# cell becomes a genuine text value (shared string) instead of a formula.
$ws.Range("A3").Formula = '="Seoul CURA Healthcare Center"'
$ws.Range("B3").Formula = '="False"'
$ws.Range("C3").Formula = '="None"'
$ws.Range("D3").Formula = '="23/12/2024"'
$ws.Range("E3").Formula = '="Some other text to fill the comment bar"'

$dataRow = $ws.Range("A3:E3")
$dataRow.Copy()
$dataRow.PasteSpecial(-4163)   # xlPasteValues

# Mirror the formatting of the row directly above (row 2) cell-by-cell via
# Paste Special -> Formats, so the new row reuses the workbook's existing
# cell styles instead of registering new ones.
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122)

$ws.Range("C2").Copy()
$ws.Range("C3").PasteSpecial(-4122)

$ws.Range("D2").Copy()
$ws.Range("D3").PasteSpecial(-4122)

$ws.Range("E2").Copy()
$ws.Range("E3").PasteSpecial(-4122)
